# Update the "ランサーズ" (Lancers) sheet: refresh the scraped-at timestamp,
# slide the job-listing window forward (rows 7-12's content moves up into
# rows 3-6, minus a couple of skill-tag columns that no longer apply), and
# drop the now-unused trailing rows 7-12 entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2026-02-02 07:00:27"

# --- Row 2: only the capture timestamp refreshes; everything else is unchanged ---
$ws.Range("A2").Value = $newTimestamp

# --- Row 3: becomes the old row 7's listing, with the refreshed timestamp ---
$ws.Range("A3").Value = $newTimestamp
$ws.Range("B3").Value = "【急募】新しいWebサービスの開発パートナーを探しています!"
$ws.Range("D3").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5483482"
$ws.Range("G3").Value = 75
$ws.Range("H3").Value = "◆開発"

# --- Row 4: becomes the old row 9's listing ---
$ws.Range("A4").Value = $newTimestamp
$ws.Range("B4").Value = "美容皮膚科向け LINE連携型BtoB SaaS(MVP) の開発案件"
$ws.Range("D4").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5483503"
$ws.Range("G4").Value = 68
$ws.Range("H4").Value = "◆開発"

# --- Row 5: new listing (Notion x Slack), and H5 no longer applies ---
$ws.Range("A5").Value = $newTimestamp
$ws.Range("B5").Value = "【急募】Notion×Slackでのオンライン講座運営システム構築"
$ws.Range("D5").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5483854"
$ws.Range("G5").Value = 28
$ws.Range("H5").ClearContents()

# --- Row 6: becomes the old row 12's listing, and H6 no longer applies ---
$ws.Range("A6").Value = $newTimestamp
$ws.Range("B6").Value = "【市場調査】海外向けデジタルサービスの価値評価依頼"
$ws.Range("D6").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5483504"
$ws.Range("G6").Value = 13
$ws.Range("H6").ClearContents()

# --- Drop the now-stale rows 7-12 (their content moved up / dropped off the window) ---
$ws.Rows("7:12").Delete()

# --- Rebuild hyperlinks for F2:F6 so the link targets match the displayed URLs ---
# (deleting any single cell's hyperlink clears the whole sheet collection in
# this host, so the only reliable way to reconcile it is delete-then-rebuild)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5483480")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5483482")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5483503")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5483854")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5483504")

# Re-adding re-applied a freshly-allocated (but equivalent) "Hyperlink" xf to
# every F cell; restore F2 to the original style index since F2 itself did
# not change in this edit.
$ws.Range("F2").Style = "Hyperlink"

# --- Narrow column B now that the longest title no longer needs as much room ---
$ws.Columns("B").ColumnWidth = 36.17
